# Updated symbol list on Sat Jan  7 22:51:18 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) values for the
# crypto rows on the active sheet, preserving their original plain-text
# cell representation (e.g. "27.10", "0.1700", "0.83%").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the refreshed feed.
$updates = @{
    "D2" = "261.34"
    "E2" = "0.83%"
    "D3" = "27.11"
    "E3" = "0.61%"
    "D4" = "4.723"
    "E4" = "0.87%"
    "D5" = "0.06205"
    "E5" = "2.61%"
    "D6" = "6.714"
    "E6" = "0.81%"
    "D7" = "0.8519"
    "E7" = "-0.96%"
    "D8" = "0.9094"
    "E8" = "-1.31%"
    "E9" = "0.76%"
    "D10" = "0.04892"
    "E10" = "-5.81%"
    "E11" = "0.04%"
    "D12" = "0.03180"
    "E12" = "3.92%"
    "D13" = "0.09061"
    "E13" = "-0.80%"
    "D14" = "0.001542"
    "E14" = "0.77%"
    "D15" = "0.0006175"
    "E15" = "1.69%"
    "D16" = "0.006144"
    "E16" = "1.22%"
    "D17" = "3.467"
    "D18" = "3.171"
    "E18" = "0.05%"
    "E19" = "-0.32%"
    "D21" = "0.1281"
    "E21" = "-1.21%"
    "D22" = "4.111"
    "E22" = "-0.03%"
    "E23" = "-0.37%"
    "D24" = "0.001219"
    "E24" = "0.14%"
    "E25" = "2.59%"
    "E26" = "0.11%"
    "D40" = "0.03912"
    "E40" = "1.47%"
    "D41" = "0.1112"
    "E41" = "-0.26%"
    "D42" = "0.004135"
    "E42" = "2.41%"
    "E43" = "-0.67%"
    "D44" = "0.01346"
    "E44" = "-10.07%"
    "D45" = "0.00005175"
    "E45" = "-0.38%"
    "E46" = "0.17%"
    "D47" = "0.03592"
    "E47" = "-34.11%"
    "D48" = "0.1698"
    "E48" = "25.48%"
    "E49" = "0.17%"
    "E50" = "0.17%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of re-interpreting it as a number/percentage, matching the
    # original inline-string cell contents (e.g. "27.10", "0.1700").
    $cell.Value = "'" + $updates[$cellRef]
    # Reset to the default "Normal" style so no quote-prefix / text
    # number-format styling gets attached to the cell (keeps styling
    # identical to the untouched cells around it).
    $cell.Style = "Normal"
}
